$p = $ppt.ActivePresentation

# 1. Update the text " & tubular necrosis " -> " and tubular necrosis " on slide 1,
#    in the "Kidney crystallopathy & tubular necrosis" rounded-rectangle shape.
$s1 = $p.Slides.Item(1)
$shp = $s1.Shapes.Item(16)
$tr = $shp.TextFrame.TextRange
$full = $tr.Text
$old = " & tubular necrosis "
$startPos = $full.IndexOf($old) + 1
$chars = $tr.Characters($startPos, $old.Length)
$chars.Text = " and tubular necrosis "

# 2. Remove the second slide from the presentation entirely.
$p.Slides.Item(2).Delete()
